$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 1.4951615
$ws.Range("H2").Value = 2.990323
$ws.Range("I2").Value = 0.6020739711267923
$ws.Range("J2").Value = 0.5021622551131893
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.6627425
$ws.Range("N2").Value = 1.325485
$ws.Range("O2").Value = 0.1766083511268686
$ws.Range("P2").Value = 0.1373859271892988
$ws.Range("Q2").Value = 0.99090707041375
$ws.Range("R2").Value = 3.963628281655
$ws.Range("S2").Value = 0.1063312912971087
$ws.Range("T2").Value = 0.06899002701819472
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 1.4951615
$ws.Range("H3").Value = 2.990323
$ws.Range("I3").Value = 0.6020739711267923
$ws.Range("J3").Value = 0.5021622551131893
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.5012916666666667
$ws.Range("N3").Value = 1.503875
$ws.Range("O3").Value = 0.1335847552912932
$ws.Range("P3").Value = 0.1558759708724027
$ws.Range("Q3").Value = 0.7495120002708334
$ws.Range("R3").Value = 4.497072001625
$ws.Range("S3").Value = 0.08042790410022965
$ws.Range("T3").Value = 0.07827502905124356
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 1.4951615
$ws.Range("H4").Value = 2.990323
$ws.Range("I4").Value = 0.6020739711267923
$ws.Range("J4").Value = 0.5021622551131893
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3493176666666667
$ws.Range("N4").Value = 1.047953
$ws.Range("O4").Value = 0.09308655643705531
$ws.Range("P4").Value = 0.1086198595652212
$ws.Range("Q4").Value = 0.5222863264698334
$ws.Range("R4").Value = 3.133717958819001
$ws.Range("S4").Value = 0.05604499269257616
$ws.Range("T4").Value = 0.05454479362934941
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 1.4951615
$ws.Range("H5").Value = 2.990323
$ws.Range("I5").Value = 0.6020739711267923
$ws.Range("J5").Value = 0.5021622551131893
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.6257723333333333
$ws.Range("N5").Value = 1.877317
$ws.Range("O5").Value = 0.166756500406739
$ws.Range("P5").Value = 0.1945830670835451
$ws.Range("Q5").Value = 0.9356307005651666
$ws.Range("R5").Value = 5.613784203391
$ws.Range("S5").Value = 0.1003997484110919
$ws.Range("T5").Value = 0.09771227177351402
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 1.4951615
$ws.Range("H6").Value = 2.990323
$ws.Range("I6").Value = 0.6020739711267923
$ws.Range("J6").Value = 0.5021622551131893
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6662906666666667
$ws.Range("N6").Value = 1.998872
$ws.Range("O6").Value = 0.1775538704869871
$ws.Range("P6").Value = 0.207182188446288
$ws.Range("Q6").Value = 0.9962121526093334
$ws.Range("R6").Value = 5.977272915656
$ws.Range("S6").Value = 0.1069005638930325
$ws.Range("T6").Value = 0.1040390749694737
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 1.4951615
$ws.Range("H7").Value = 2.990323
$ws.Range("I7").Value = 0.6020739711267923
$ws.Range("J7").Value = 0.5021622551131893
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.9471965
$ws.Range("N7").Value = 1.894393
$ws.Range("O7").Value = 0.2524099662510568
$ws.Range("P7").Value = 0.196352986843244
$ws.Range("Q7").Value = 1.41621173973475
$ws.Range("R7").Value = 5.664846958939
$ws.Range("S7").Value = 0.1519694707327534
$ws.Range("T7").Value = 0.0986010586714138
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.9881903333333333
$ws.Range("H8").Value = 2.964571
$ws.Range("I8").Value = 0.3979260288732077
$ws.Range("J8").Value = 0.4978377448868108
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 0.6627425
$ws.Range("N8").Value = 1.325485
$ws.Range("O8").Value = 0.1766083511268686
$ws.Range("P8").Value = 0.1373859271892988
$ws.Range("Q8").Value = 0.6549157319891666
$ws.Range("R8").Value = 3.929494391935
$ws.Range("S8").Value = 0.07027705982975993
$ws.Range("T8").Value = 0.06839590017110408
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.9881903333333333
$ws.Range("H9").Value = 2.964571
$ws.Range("I9").Value = 0.3979260288732077
$ws.Range("J9").Value = 0.4978377448868108
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5012916666666667
$ws.Range("N9").Value = 1.503875
$ws.Range("O9").Value = 0.1335847552912932
$ws.Range("P9").Value = 0.1558759708724027
$ws.Range("Q9").Value = 0.4953715791805556
$ws.Range("R9").Value = 4.458344212625
$ws.Range("S9").Value = 0.05315685119106351
$ws.Range("T9").Value = 0.07760094182115916
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.9881903333333333
$ws.Range("H10").Value = 2.964571
$ws.Range("I10").Value = 0.3979260288732077
$ws.Range("J10").Value = 0.4978377448868108
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3493176666666667
$ws.Range("N10").Value = 1.047953
$ws.Range("O10").Value = 0.09308655643705531
$ws.Range("P10").Value = 0.1086198595652212
$ws.Range("Q10").Value = 0.3451923414625556
$ws.Range("R10").Value = 3.106731073163
$ws.Range("S10").Value = 0.03704156374447915
$ws.Range("T10").Value = 0.05407506593587181
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.9881903333333333
$ws.Range("H11").Value = 2.964571
$ws.Range("I11").Value = 0.3979260288732077
$ws.Range("J11").Value = 0.4978377448868108
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.6257723333333333
$ws.Range("N11").Value = 1.877317
$ws.Range("O11").Value = 0.166756500406739
$ws.Range("P11").Value = 0.1945830670835451
$ws.Range("Q11").Value = 0.6183821706674444
$ws.Range("R11").Value = 5.565439536006999
$ws.Range("S11").Value = 0.06635675199564708
$ws.Range("T11").Value = 0.09687079531003112
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.9881903333333333
$ws.Range("H12").Value = 2.964571
$ws.Range("I12").Value = 0.3979260288732077
$ws.Range("J12").Value = 0.4978377448868108
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.6662906666666667
$ws.Range("N12").Value = 1.998872
$ws.Range("O12").Value = 0.1775538704869871
$ws.Range("P12").Value = 0.207182188446288
$ws.Range("Q12").Value = 0.6584219959902222
$ws.Range("R12").Value = 5.925797963911999
$ws.Range("S12").Value = 0.07065330659395462
$ws.Range("T12").Value = 0.1031431134768143
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.9881903333333333
$ws.Range("H13").Value = 2.964571
$ws.Range("I13").Value = 0.3979260288732077
$ws.Range("J13").Value = 0.4978377448868108
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 0.9471965
$ws.Range("N13").Value = 1.894393
$ws.Range("O13").Value = 0.2524099662510568
$ws.Range("P13").Value = 0.196352986843244
$ws.Range("Q13").Value = 0.9360104250671666
$ws.Range("R13").Value = 5.616062550403
$ws.Range("S13").Value = 0.1004404955183034
$ws.Range("T13").Value = 0.09775192817183022
